$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

$emuPerPoint = 12700
$left   = 5569075 / $emuPerPoint
$top    = 6031150 / $emuPerPoint
$width  = 3400200 / $emuPerPoint
$height = 364800 / $emuPerPoint
$inset  = 91425 / $emuPerPoint

$shape = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shape.Name = "Shape 213"

$shape.Fill.Visible = 0
$shape.Line.Visible = 0

$tf = $shape.TextFrame
$tf.AutoSize = 0
$tf.VerticalAnchor = 1
$tf.HorizontalAnchor = 0
$tf.MarginLeft = $inset
$tf.MarginRight = $inset
$tf.MarginTop = $inset
$tf.MarginBottom = $inset

$tr = $tf.TextRange
$tr.Text = "https://github.com/EddyTheB/HarmonicSeries"
$tr.Font.Size = 12
$tr.ParagraphFormat.SpaceBefore = 0
$tr.ParagraphFormat.Bullet.Visible = 0
$tr.IndentLevel = 0
